$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "nironi1"
$ws.Range("B11").Value = "nir12345@"
$ws.Range("C11").Value = 209375900
